$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Column -> new value
$updates = @{
    2  = @{ E=3; G=36.899643; H=110.698929; I=0.7238945645409351; J=0.7238945645409351; K=3;
            M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588;
            Q=8655.969262546507; R=77903.72336291856; S=0.7081676734686511; T=0.7081676734686511 }
    3  = @{ E=3; G=36.899643; H=110.698929; I=0.7238945645409351; J=0.7238945645409351; K=3;
            M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783;
            Q=77.87054661099999; R=700.8349194989999; S=0.006370794783647458; T=0.006370794783647458 }
    4  = @{ E=3; G=36.899643; H=110.698929; I=0.7238945645409351; J=0.7238945645409351; K=3;
            M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245;
            Q=114.360037779174; R=1029.240340012566; S=0.009356096288636696; T=0.009356096288636696 }
    5  = @{ E=3; G=3.374819; H=10.124457; I=0.0662069584361419; J=0.0662069584361419; K=3;
            M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588;
            Q=791.6697061447976; R=7125.027355303178; S=0.06476858650388025; T=0.06476858650388025 }
    6  = @{ E=3; G=3.374819; H=10.124457; I=0.0662069584361419; J=0.0662069584361419; K=3;
            M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783;
            Q=7.121993029666666; R=64.09793726699999; S=0.0005826690323522732; T=0.0005826690323522732 }
    7  = @{ E=3; G=3.374819; H=10.124457; I=0.0662069584361419; J=0.0662069584361419; K=3;
            M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245;
            Q=10.459299791542; R=94.13369812387799; S=0.0008557028999093732; T=0.0008557028999093732 }
    8  = @{ E=3; G=10.699319; H=32.097957; I=0.2098984770229228; J=0.2098984770229228; K=3;
            M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588;
            Q=2509.861041045298; R=22588.74936940768; S=0.2053383509409274; T=0.2053383509409274 }
    9  = @{ E=3; G=10.699319; H=32.097957; I=0.2098984770229228; J=0.2098984770229228; K=3;
            M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783;
            Q=22.57912952966667; R=203.212165767; S=0.001847258134009051; T=0.001847258134009051 }
    10 = @{ E=3; G=10.699319; H=32.097957; I=0.2098984770229228; J=0.2098984770229228; K=3;
            M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245;
            Q=33.159522032542; R=298.435698292878; S=0.002712867947986383; T=0.002712867947986383 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
